$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.034663131525051
$ws.Cells.Item(2, 4).Value = 1.04221534779115
$ws.Cells.Item(2, 5).Value = 1.033813452028654
$ws.Cells.Item(2, 6).Value = 1.050712351946554
$ws.Cells.Item(2, 9).Value = 1.034566327672007
$ws.Cells.Item(2, 10).Value = 1.039781144244546
$ws.Cells.Item(2, 11).Value = 1.044992664679298
$ws.Cells.Item(2, 12).Value = 1.036614721484526
$ws.Cells.Item(2, 13).Value = 1.053465867625742
$ws.Cells.Item(2, 14).Value = 1.041257752733938
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.035611391369041
$ws.Cells.Item(3, 4).Value = 1.043064763084082
$ws.Cells.Item(3, 5).Value = 1.034619404088002
$ws.Cells.Item(3, 6).Value = 1.05165565324506
$ws.Cells.Item(3, 9).Value = 1.034692445792059
$ws.Cells.Item(3, 10).Value = 1.0403726357229
$ws.Cells.Item(3, 11).Value = 1.045652764882827
$ws.Cells.Item(3, 12).Value = 1.037229740026358
$ws.Cells.Item(3, 13).Value = 1.054221333256409
$ws.Cells.Item(3, 14).Value = 1.041850084198036
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.036225446064468
$ws.Cells.Item(4, 4).Value = 1.043615124112124
$ws.Cells.Item(4, 5).Value = 1.035141682476418
$ws.Cells.Item(4, 6).Value = 1.0522668618481
$ws.Cells.Item(4, 9).Value = 1.034772717159129
$ws.Cells.Item(4, 10).Value = 1.040755223629329
$ws.Cells.Item(4, 11).Value = 1.046079982262542
$ws.Cells.Item(4, 12).Value = 1.037627819711004
$ws.Cells.Item(4, 13).Value = 1.054710373592719
$ws.Cells.Item(4, 14).Value = 1.042233215423176
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.036483705468928
$ws.Cells.Item(5, 4).Value = 1.043846669868189
$ws.Cells.Item(5, 5).Value = 1.035361432086817
$ws.Cells.Item(5, 6).Value = 1.05252401112209
$ws.Cells.Item(5, 9).Value = 1.034806142987748
$ws.Cells.Item(5, 10).Value = 1.040916027455054
$ws.Cells.Item(5, 11).Value = 1.046259604489462
$ws.Cells.Item(5, 12).Value = 1.037795200421467
$ws.Cells.Item(5, 13).Value = 1.054916013630361
$ws.Cells.Item(5, 14).Value = 1.042394247608777
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.036527074827453
$ws.Cells.Item(6, 4).Value = 1.043885557581338
$ws.Cells.Item(6, 5).Value = 1.035398339749213
$ws.Cells.Item(6, 6).Value = 1.052567199137815
$ws.Cells.Item(6, 9).Value = 1.03481173653606
$ws.Cells.Item(6, 10).Value = 1.040943024976461
$ws.Cells.Item(6, 11).Value = 1.046289764994242
$ws.Cells.Item(6, 12).Value = 1.037823305990124
$ws.Cells.Item(6, 13).Value = 1.054950544245444
$ws.Cells.Item(6, 14).Value = 1.042421283469762
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.036228896505142
$ws.Cells.Item(7, 4).Value = 1.043618217355971
$ws.Cells.Item(7, 5).Value = 1.035144618061117
$ws.Cells.Item(7, 6).Value = 1.052270297116302
$ws.Cells.Item(7, 9).Value = 1.034773165055733
$ws.Cells.Item(7, 10).Value = 1.040757372441171
$ws.Cells.Item(7, 11).Value = 1.046082382306194
$ws.Cells.Item(7, 12).Value = 1.037630056151554
$ws.Cells.Item(7, 13).Value = 1.054713121179971
$ws.Cells.Item(7, 14).Value = 1.042235367286577
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.034983503028142
$ws.Cells.Item(8, 4).Value = 1.042502259154133
$ws.Cells.Item(8, 5).Value = 1.034085666294974
$ws.Cells.Item(8, 6).Value = 1.051030972407069
$ws.Cells.Item(8, 9).Value = 1.034609226081254
$ws.Cells.Item(8, 10).Value = 1.039981071303007
$ws.Cells.Item(8, 11).Value = 1.045215729651785
$ws.Cells.Item(8, 12).Value = 1.036822544053629
$ws.Cells.Item(8, 13).Value = 1.053721137872151
$ws.Cells.Item(8, 14).Value = 1.041457963711754
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.032792584059919
$ws.Cells.Item(9, 4).Value = 1.040541469944149
$ws.Cells.Item(9, 5).Value = 1.032225644186658
$ws.Cells.Item(9, 6).Value = 1.048853540596575
$ws.Cells.Item(9, 9).Value = 1.034310140037684
$ws.Cells.Item(9, 10).Value = 1.038612051772678
$ws.Cells.Item(9, 11).Value = 1.04368930368865
$ws.Cells.Item(9, 12).Value = 1.035400583757839
$ws.Cells.Item(9, 13).Value = 1.051974750801814
$ws.Cells.Item(9, 14).Value = 1.040087000016663
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.031334467729799
$ws.Cells.Item(10, 4).Value = 1.039238172546806
$ws.Cells.Item(10, 5).Value = 1.030989737259377
$ws.Cells.Item(10, 6).Value = 1.047406317752893
$ws.Cells.Item(10, 9).Value = 1.034103919935111
$ws.Cells.Item(10, 10).Value = 1.037698707879896
$ws.Cells.Item(10, 11).Value = 1.042672241268603
$ws.Cells.Item(10, 12).Value = 1.034453332943965
$ws.Cells.Item(10, 13).Value = 1.050811651863568
$ws.Cells.Item(10, 14).Value = 1.039172359070791
$ws.Cells.Item(11, 2).Value = 1.019999999999999
$ws.Cells.Item(11, 3).Value = 1.030703691617573
$ws.Cells.Item(11, 4).Value = 1.038674771410601
$ws.Cells.Item(11, 5).Value = 1.030455567562953
$ws.Cells.Item(11, 6).Value = 1.04678071500069
$ws.Cells.Item(11, 9).Value = 1.0340130110744
$ws.Cells.Item(11, 10).Value = 1.037303074873275
$ws.Cells.Item(11, 11).Value = 1.042231987383964
$ws.Cells.Item(11, 12).Value = 1.034043347688162
$ws.Cells.Item(11, 13).Value = 1.050308307426516
$ws.Cells.Item(11, 14).Value = 1.03877616421992
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.030469484045189
$ws.Cells.Item(12, 4).Value = 1.038465641006348
$ws.Cells.Item(12, 5).Value = 1.030257302692979
$ws.Cells.Item(12, 6).Value = 1.046548498059354
$ws.Cells.Item(12, 9).Value = 1.033979001461811
$ws.Cells.Item(12, 10).Value = 1.037156097740705
$ws.Cells.Item(12, 11).Value = 1.042068479731042
$ws.Cells.Item(12, 12).Value = 1.03389108909678
$ws.Cells.Item(12, 13).Value = 1.050121386924622
$ws.Cells.Item(12, 14).Value = 1.038628978362964
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.030519718242295
$ws.Cells.Item(13, 4).Value = 1.038510493746613
$ws.Cells.Item(13, 5).Value = 1.030299824391499
$ws.Cells.Item(13, 6).Value = 1.046598302129172
$ws.Cells.Item(13, 9).Value = 1.033986307582756
$ws.Cells.Item(13, 10).Value = 1.037187625788074
$ws.Cells.Item(13, 11).Value = 1.042103551653776
$ws.Cells.Item(13, 12).Value = 1.03392374778902
$ws.Cells.Item(13, 13).Value = 1.05016147999485
$ws.Cells.Item(13, 14).Value = 1.038660551183777
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.03068433008965
$ws.Cells.Item(14, 4).Value = 1.038657481724347
$ws.Cells.Item(14, 5).Value = 1.030439175868419
$ws.Cells.Item(14, 6).Value = 1.04676151660851
$ws.Cells.Item(14, 9).Value = 1.034010204763328
$ws.Cells.Item(14, 10).Value = 1.037290926128913
$ws.Cells.Item(14, 11).Value = 1.042218471330383
$ws.Cells.Item(14, 12).Value = 1.034030761364391
$ws.Cells.Item(14, 13).Value = 1.050292855614285
$ws.Cells.Item(14, 14).Value = 1.038763998222948
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.030785764936484
$ws.Cells.Item(15, 4).Value = 1.038748064707411
$ws.Cells.Item(15, 5).Value = 1.030525054769194
$ws.Cells.Item(15, 6).Value = 1.04686209964916
$ws.Cells.Item(15, 9).Value = 1.034024896550774
$ws.Cells.Item(15, 10).Value = 1.037354570069388
$ws.Cells.Item(15, 11).Value = 1.042289280115598
$ws.Cells.Item(15, 12).Value = 1.034096699726699
$ws.Cells.Item(15, 13).Value = 1.050373806340245
$ws.Cells.Item(15, 14).Value = 1.038827732545119
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.031376342939419
$ws.Cells.Item(16, 4).Value = 1.03927558344467
$ws.Cells.Item(16, 5).Value = 1.031025209237919
$ws.Cells.Item(16, 6).Value = 1.047447859277671
$ws.Cells.Item(16, 9).Value = 1.034109919290261
$ws.Cells.Item(16, 10).Value = 1.037724961682231
$ws.Cells.Item(16, 11).Value = 1.042701462528829
$ws.Cells.Item(16, 12).Value = 1.034480546210123
$ws.Cells.Item(16, 13).Value = 1.050845063285118
$ws.Cells.Item(16, 14).Value = 1.039198650156536
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.031746957646924
$ws.Cells.Item(17, 4).Value = 1.039606733357253
$ws.Cells.Item(17, 5).Value = 1.031339207919289
$ws.Cells.Item(17, 6).Value = 1.047815574054535
$ws.Cells.Item(17, 9).Value = 1.034162819995601
$ws.Cells.Item(17, 10).Value = 1.037957259334395
$ws.Cells.Item(17, 11).Value = 1.04296005214016
$ws.Cells.Item(17, 12).Value = 1.034721372063492
$ws.Cells.Item(17, 13).Value = 1.051140747489123
$ws.Cells.Item(17, 14).Value = 1.039431277698012
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.031963188559301
$ws.Cells.Item(18, 4).Value = 1.039799977644641
$ws.Cells.Item(18, 5).Value = 1.031522453095201
$ws.Cells.Item(18, 6).Value = 1.048030157493589
$ws.Cells.Item(18, 9).Value = 1.034193520301934
$ws.Cells.Item(18, 10).Value = 1.038092740175603
$ws.Cells.Item(18, 11).Value = 1.043110896603412
$ws.Cells.Item(18, 12).Value = 1.034861859063616
$ws.Cells.Item(18, 13).Value = 1.051313242507888
$ws.Cells.Item(18, 14).Value = 1.039566950937554
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.032036927458098
$ws.Cells.Item(19, 4).Value = 1.039865884201888
$ws.Cells.Item(19, 5).Value = 1.031584951046121
$ws.Cells.Item(19, 6).Value = 1.048103342097363
$ws.Cells.Item(19, 9).Value = 1.034203961874807
$ws.Cells.Item(19, 10).Value = 1.038138933160746
$ws.Cells.Item(19, 11).Value = 1.043162332927238
$ws.Cells.Item(19, 12).Value = 1.034909764411666
$ws.Cells.Item(19, 13).Value = 1.051372063466801
$ws.Cells.Item(19, 14).Value = 1.039613209522034
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.03170718822917
$ws.Cells.Item(20, 4).Value = 1.039571194771579
$ws.Cells.Item(20, 5).Value = 1.031305508990204
$ws.Cells.Item(20, 6).Value = 1.047776111216831
$ws.Cells.Item(20, 9).Value = 1.034157160360008
$ws.Cells.Item(20, 10).Value = 1.037932337482215
$ws.Cells.Item(20, 11).Value = 1.042932306516082
$ws.Cells.Item(20, 12).Value = 1.034695531925907
$ws.Cells.Item(20, 13).Value = 1.051109020526931
$ws.Cells.Item(20, 14).Value = 1.039406320453942
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.030635853490457
$ws.Cells.Item(21, 4).Value = 1.038614193501487
$ws.Cells.Item(21, 5).Value = 1.030398136192617
$ws.Cells.Item(21, 6).Value = 1.046713449590669
$ws.Cells.Item(21, 9).Value = 1.034003174316017
$ws.Cells.Item(21, 10).Value = 1.037260507336619
$ws.Cells.Item(21, 11).Value = 1.042184629727249
$ws.Cells.Item(21, 12).Value = 1.033999247747112
$ws.Cells.Item(21, 13).Value = 1.050254167543602
$ws.Cells.Item(21, 14).Value = 1.038733536232481
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.029962788401241
$ws.Cells.Item(22, 4).Value = 1.038013310267042
$ws.Cells.Item(22, 5).Value = 1.029828500771716
$ws.Cells.Item(22, 6).Value = 1.046046237147549
$ws.Cells.Item(22, 9).Value = 1.033904957208665
$ws.Cells.Item(22, 10).Value = 1.036837977624305
$ws.Cells.Item(22, 11).Value = 1.041714664876495
$ws.Cells.Item(22, 12).Value = 1.033561629853286
$ws.Cells.Item(22, 13).Value = 1.049716942904658
$ws.Cells.Item(22, 14).Value = 1.03831040647951
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.03031954280687
$ws.Cells.Item(23, 4).Value = 1.038331771522642
$ws.Cells.Item(23, 5).Value = 1.030130392750985
$ws.Cells.Item(23, 6).Value = 1.046399850937642
$ws.Cells.Item(23, 9).Value = 1.033957156489127
$ws.Cells.Item(23, 10).Value = 1.037061979992598
$ws.Cells.Item(23, 11).Value = 1.041963789526242
$ws.Cells.Item(23, 12).Value = 1.033793603518763
$ws.Cells.Item(23, 13).Value = 1.050001711242904
$ws.Cells.Item(23, 14).Value = 1.038534726956859
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.031725158139201
$ws.Cells.Item(24, 4).Value = 1.039587252850301
$ws.Cells.Item(24, 5).Value = 1.031320735792135
$ws.Cells.Item(24, 6).Value = 1.047793942460637
$ws.Cells.Item(24, 9).Value = 1.034159718187097
$ws.Cells.Item(24, 10).Value = 1.037943598639196
$ws.Cells.Item(24, 11).Value = 1.042944843527635
$ws.Cells.Item(24, 12).Value = 1.034707207918283
$ws.Cells.Item(24, 13).Value = 1.051123356490843
$ws.Cells.Item(24, 14).Value = 1.039417597603058
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.033358553553026
$ws.Cells.Item(25, 4).Value = 1.041047700371609
$ws.Cells.Item(25, 5).Value = 1.032705786533502
$ws.Cells.Item(25, 6).Value = 1.049415689837300
$ws.Cells.Item(25, 9).Value = 1.034388666996105
$ws.Cells.Item(25, 10).Value = 1.038966097632703
$ws.Cells.Item(25, 11).Value = 1.044083828439635
$ws.Cells.Item(25, 12).Value = 1.035768071907712
$ws.Cells.Item(25, 13).Value = 1.052426034824406
$ws.Cells.Item(25, 14).Value = 1.040441548662419
